$wb = $excel.ActiveWorkbook

# Update values in the "dEX4_1" worksheet (column B, rows 2-13)
$ws2 = $wb.Worksheets.Item("dEX4_1")

$ws2.Range("B2").Value = 24102
$ws2.Range("B3").Value = 34545
$ws2.Range("B4").Value = 658817
$ws2.Range("B5").Value = 10613437
$ws2.Range("B6").Value = 7169087
$ws2.Range("B7").Value = 5161263
$ws2.Range("B8").Value = 7985127
$ws2.Range("B9").Value = 7405939
$ws2.Range("B10").Value = 1476866
$ws2.Range("B11").Value = 557280
$ws2.Range("B12").Value = 70910
$ws2.Range("B13").Value = 28789

# Make "dEX4_1" the active/selected sheet (this also clears tabSelected on
# whichever sheet previously had it, and updates the workbook's activeTab)
$ws2.Activate()

# Select E21:E22 on the now-active sheet (matches the saved selection state)
$ws2.Range("E21:E22").Select()

$wb.Save()
